$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_data")

# Update cell C10 value from 1010 to 6010
$ws.Range("C10").Value = 6010

# Clear the active selection (set selection back to A1 so no <selection> override is written)
$ws.Range("A1").Select()
